$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K6").Value = "95b8b0b0-a12d-48c7-ac7d-b79791c1de3c"
$ws.Range("K7").Value = "2635d729-6559-464a-b0d6-7325018ffe36"
$ws.Range("K8").Value = "b9069bd0-8299-4b62-ab81-e287d88133be"
$ws.Range("K9").Value = "ffd88615-fc9c-402f-8479-67df9fcff002"
$ws.Range("K10").Value = "23166d7b-d6b5-4bf5-906e-c0350d65bcf1"
$ws.Range("K11").Value = "6f584c45-31fb-419e-9173-33dee34cffa9"
$ws.Range("K12").Value = "6d58188e-2b9a-42db-b593-b578e796f775"
$ws.Range("K13").Value = "ef252f02-7ff7-4fd4-9c14-f9e2508c6ffa"
$ws.Range("K14").Value = "cfeeb1dd-de6e-423c-9f9d-ce724e1979dd"
$ws.Range("K15").Value = "499b7f18-15d5-405a-9a99-830fcf1a67e7"
$ws.Range("K16").Value = "57858b87-1b4d-4bec-870b-7bfe9bcbf4c6"
$ws.Range("K17").Value = "fda5965b-8733-4cd7-8f38-46707a87000e"
$ws.Range("K19").Value = "bbf9588f-c3b5-490f-8e1d-0b4d585a350a"
$ws.Range("K20").Value = "3645688c-9317-4f40-b626-c4f2f55b9b2b"
$ws.Range("K21").Value = "7b6fcfa1-dab6-4c98-b71a-e2dcea765514"
$ws.Range("K22").Value = "9686114b-e3f3-4f48-a423-5fb7cfd98da6"
$ws.Range("K23").Value = "b11d057d-511d-4e08-b2a0-7956484b2bb4"
$ws.Range("K24").Value = "b72d9b61-880e-48f9-b52a-3e2e80071be7"
$ws.Range("K25").Value = "23b92792-6856-4dd6-8f36-26291e6f7c3b"
$ws.Range("K26").Value = "0f74e96e-82c9-4553-b393-f514c76dcb21"
$ws.Range("K27").Value = "962da885-0a5f-40d9-bd3c-9787ed79257f"
$ws.Range("K28").Value = "7647913f-db03-4601-af23-12d6bb2410dc"
$ws.Range("K29").Value = "bb824194-3bd3-45d0-bd6b-673b34cf993b"
$ws.Range("K30").Value = "2599baf6-9da4-4128-8a54-934e6cadee6e"
$ws.Range("K31").Value = "05ae4641-e06a-4d39-9ebd-0935ad2e061b"
$ws.Range("K32").Value = "9bc6e87f-e77f-4750-a7c9-b1b320902717"
$ws.Range("K33").Value = "4cf83cfc-cca5-4f4f-8c6e-f0e901ec1e67"
$ws.Range("K34").Value = "e773ff1a-56fe-4d8f-99e4-cfd83382a9e6"
$ws.Range("K35").Value = "04a1b1b9-2e96-47c3-9ca4-f53d9194c2d1"
$ws.Range("K36").Value = "5c2ca29d-4b10-47c0-8212-3150a077625a"
$ws.Range("K37").Value = "dc5f13d6-c58f-440c-8b16-5187f61f1dbc"
$ws.Range("K38").Value = "846e6fd6-9a03-4e19-8866-bcc60e15178a"
$ws.Range("K39").Value = "54f0d4e0-2a53-49c0-9871-ac29f241f5ff"
$ws.Range("K40").Value = "66a4e63a-4cd2-47ff-8d31-fd14175a5169"
$ws.Range("K41").Value = "c580c9af-cce5-4c1e-94aa-54625dd203e3"
$ws.Range("K42").Value = "88bec9e8-b860-4b77-a028-80c29e53e247"
$ws.Range("K43").Value = "edc93ddc-555f-4661-881e-b1bdbfde4d3a"
$ws.Range("K44").Value = "af5f9e17-b5c1-4945-8069-ac7535280835"
$ws.Range("K45").Value = "9fbea7ac-4085-4f12-8aae-d5f6212aaf6c"
$ws.Range("K46").Value = "8255fed7-726b-4856-b95f-54780f6f26e7"
$ws.Range("K47").Value = "6fa663c9-6616-41b2-9cb3-8e04ad299ac7"
$ws.Range("K48").Value = "44af2d82-e3e2-4fbf-b43d-65a7b4863bf3"
$ws.Range("K49").Value = "ae8b4b12-bc48-4f77-b3ba-bcdf37ed56ed"
$ws.Range("K50").Value = "4edd6a35-2dd9-4580-ab45-9e2d7b4c962a"
$ws.Range("K51").Value = "777c5e3b-469a-4d0a-baf5-76b83f79c3d4"
$ws.Range("K52").Value = "5e5a732a-971e-41fd-9ec2-5ea31c167e7d"
$ws.Range("K53").Value = "f180b0a2-8968-4324-899c-81d7805dc46b"
$ws.Range("K54").Value = "a0dbf012-0e14-46b4-b284-835c7fde1d83"
$ws.Range("K55").Value = "5556aa4f-0093-4423-b16b-ba324c3d003e"
$ws.Range("K56").Value = "805e628a-64e0-4a8e-8a86-83ecd5d662cf"
$ws.Range("K57").Value = "81338169-7388-4b32-911f-e552e099f102"
$ws.Range("K58").Value = "f718546c-5f2a-48b7-a86c-fb7f9aaced33"
$ws.Range("K59").Value = "361483ad-3fd1-4c32-8ab0-e6e8a2bb3b1a"
$ws.Range("K60").Value = "924d09c8-a080-4da4-b566-4d1884148ec1"
$ws.Range("K61").Value = "b859a5ff-f24c-4616-b10b-ee6a97d935d6"
$ws.Range("K62").Value = "793aa431-f23c-48be-b929-50e1e2c164b1"
$ws.Range("K63").Value = "6db8e041-81f3-4ddc-9b21-de3999a7f22f"
$ws.Range("K64").Value = "8f3e0493-cdd9-4d9d-8a80-e85f36908d3b"
$ws.Range("K65").Value = "70fc27e8-4c10-4045-858c-b294df3bde81"
$ws.Range("K66").Value = "89be468e-3403-4958-a2cf-0284d7287d67"
